$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 3, 4, 6, 7, 8
$ws.Range("F3").Value = -6
$ws.Range("F4").Value = -7
$ws.Range("F6").Value = -5
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 0
